$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.49"
$ws.Range("E2").Value = "'0.47%"
$ws.Range("D3").Value = "'37.20"
$ws.Range("E3").Value = "'-1.79%"
$ws.Range("D4").Value = "'5.132"
$ws.Range("E4").Value = "'1.70%"
$ws.Range("D5").Value = "'0.07767"
$ws.Range("E5").Value = "'-1.63%"
$ws.Range("E6").Value = "'-0.16%"
$ws.Range("D7").Value = "'1.876"
$ws.Range("E7").Value = "'-7.30%"
$ws.Range("D8").Value = "'0.9179"
$ws.Range("E8").Value = "'-0.77%"
$ws.Range("D9").Value = "'0.1189"
$ws.Range("E9").Value = "'-6.63%"
$ws.Range("D10").Value = "'0.1886"
$ws.Range("E10").Value = "'-0.77%"
$ws.Range("D11").Value = "'0.09547"
$ws.Range("E11").Value = "'9.67%"
$ws.Range("D12").Value = "'0.03438"
$ws.Range("E12").Value = "'-0.62%"
$ws.Range("D13").Value = "'0.09692"
$ws.Range("E13").Value = "'-0.35%"
$ws.Range("D14").Value = "'0.001373"
$ws.Range("E14").Value = "'-1.42%"
$ws.Range("D15").Value = "'0.005782"
$ws.Range("E15").Value = "'-2.70%"
$ws.Range("D16").Value = "'3.539"
$ws.Range("E16").Value = "'-0.21%"
$ws.Range("D17").Value = "'4.401"
$ws.Range("E17").Value = "'0.34%"
$ws.Range("D18").Value = "'3.050"
$ws.Range("E18").Value = "'-4.00%"
$ws.Range("E19").Value = "'-1.05%"
$ws.Range("D20").Value = "'5.256"
$ws.Range("E20").Value = "'4.92%"
$ws.Range("D21").Value = "'0.1267"
$ws.Range("E21").Value = "'-2.46%"
$ws.Range("E22").Value = "'3.12%"
$ws.Range("D23").Value = "'0.02107"
$ws.Range("E23").Value = "'5,598.21%"
$ws.Range("D24").Value = "'0.04336"
$ws.Range("E24").Value = "'0.19%"
$ws.Range("D25").Value = "'0.001199"
$ws.Range("E25").Value = "'-2.07%"
$ws.Range("E26").Value = "'-7.47%"
$ws.Range("E27").Value = "'-63.75%"
$ws.Range("D39").Value = "'0.02062"
$ws.Range("E39").Value = "'-8.52%"
$ws.Range("D40").Value = "'0.05013"
$ws.Range("E40").Value = "'0.21%"
$ws.Range("D41").Value = "'0.007684"
$ws.Range("E41").Value = "'1.38%"
$ws.Range("D42").Value = "'0.009822"
$ws.Range("E42").Value = "'-0.51%"
$ws.Range("D43").Value = "'0.1342"
$ws.Range("E43").Value = "'-1.01%"
$ws.Range("D44").Value = "'0.002174"
$ws.Range("E44").Value = "'3.84%"
$ws.Range("D45").Value = "'0.008760"
$ws.Range("E45").Value = "'2.58%"
$ws.Range("D46").Value = "'0.00006710"
$ws.Range("E46").Value = "'4.59%"
$ws.Range("E47").Value = "'0.00%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.002938"
$ws.Range("E48").Value = "'-2.14%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.001201"
$ws.Range("E49").Value = "'-0.05%"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("E51").Value = "'0.00%"
